{"js": "// Update the date heading.\nconst dateResults = context.document.body.search(\"2024-12-23 Monday\", {\n  matchCase: true,\n  matchWildcards: false,\n});\ndateResults.load(\"items\");\nawait context.sync();\ndateResults.items[0].insertText(\"2024-12-24 Tuesday\", Word.InsertLocation.replace);\n\n// Update the practice-problem table. Cell positions are addressed by\n// (tableRow, column) so the duplicate \"20\u00f73=\" values in the sheet\n// (row 8 col 2 and row 16 col 4) are disambiguated correctly.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst cellUpdates = [\n  // row 0\n  [0, 0, \"79\u00f77=\"],\n  [0, 1, \"68\u00f73=\"],\n  [0, 2, \"97\u00f76=\"],\n  [0, 3, \"25\u00f78=\"],\n  [0, 4, \"48\u00f77=\"],\n  // row 4\n  [4, 0, \"53\u00f74=\"],\n  [4, 1, \"69\u00f73=\"],\n  [4, 2, \"54\u00f77=\"],\n  [4, 3, \"65\u00f77=\"],\n  [4, 4, \"62\u00f77=\"],\n  // row 8\n  [8, 0, \"82\u00f76=\"],\n  [8, 1, \"84\u00f73=\"],\n  [8, 2, \"32\u00f78=\"],\n  [8, 3, \"41\u00f72=\"],\n  [8, 4, \"48\u00f72=\"],\n  // row 12\n  [12, 0, \"83\u00f73=\"],\n  [12, 1, \"49\u00f74=\"],\n  [12, 2, \"92\u00f79=\"],\n  [12, 3, \"62\u00f79=\"],\n  [12, 4, \"93\u00f79=\"],\n  // row 16\n  [16, 0, \"89\u00f73=\"],\n  [16, 1, \"98\u00f78=\"],\n  [16, 2, \"51\u00f72=\"],\n  [16, 3, \"21\u00f78=\"],\n  [16, 4, \"27\u00f73=\"],\n];\n\n// First load the current text of every target cell so we can build a\n// search query that is guaranteed to match (and only match) that cell's\n// existing run, which keeps the run's formatting (font/size/alignment)\n// intact when it is replaced.\nconst cells = cellUpdates.map(([row, col]) => table.getCell(row, col));\nfor (const cell of cells) {\n  cell.body.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < cellUpdates.length; i++) {\n  const [, , newText] = cellUpdates[i];\n  const cell = cells[i];\n  const oldText = cell.body.text;\n  const cellResults = cell.body.search(oldText, {\n    matchCase: true,\n    matchWildcards: false,\n  });\n  cellResults.load(\"items\");\n  await context.sync();\n  cellResults.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2024-12-24 Tuesday\"\n\n# Update the practice-problem table. Table.Cell(row, col) is 1-indexed;\n# the five data rows of the 20-row table (blank spacer rows in between)\n# are rows 1, 5, 9, 13, 17, each with 5 columns.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text  = \"79\u00f77=\"\n$t.Cell(1, 2).Range.Text  = \"68\u00f73=\"\n$t.Cell(1, 3).Range.Text  = \"97\u00f76=\"\n$t.Cell(1, 4).Range.Text  = \"25\u00f78=\"\n$t.Cell(1, 5).Range.Text  = \"48\u00f77=\"\n\n$t.Cell(5, 1).Range.Text  = \"53\u00f74=\"\n$t.Cell(5, 2).Range.Text  = \"69\u00f73=\"\n$t.Cell(5, 3).Range.Text  = \"54\u00f77=\"\n$t.Cell(5, 4).Range.Text  = \"65\u00f77=\"\n$t.Cell(5, 5).Range.Text  = \"62\u00f77=\"\n\n$t.Cell(9, 1).Range.Text  = \"82\u00f76=\"\n$t.Cell(9, 2).Range.Text  = \"84\u00f73=\"\n$t.Cell(9, 3).Range.Text  = \"32\u00f78=\"\n$t.Cell(9, 4).Range.Text  = \"41\u00f72=\"\n$t.Cell(9, 5).Range.Text  = \"48\u00f72=\"\n\n$t.Cell(13, 1).Range.Text = \"83\u00f73=\"\n$t.Cell(13, 2).Range.Text = \"49\u00f74=\"\n$t.Cell(13, 3).Range.Text = \"92\u00f79=\"\n$t.Cell(13, 4).Range.Text = \"62\u00f79=\"\n$t.Cell(13, 5).Range.Text = \"93\u00f79=\"\n\n$t.Cell(17, 1).Range.Text = \"89\u00f73=\"\n$t.Cell(17, 2).Range.Text = \"98\u00f78=\"\n$t.Cell(17, 3).Range.Text = \"51\u00f72=\"\n$t.Cell(17, 4).Range.Text = \"21\u00f78=\"\n$t.Cell(17, 5).Range.Text = \"27\u00f73=\"\n"}
